# "Actualizando hitos correctos Febrero"
# Adds a new milestone row ("Febrero" / 27-02-15 / 27-02-15) to the
# "Hitos" table on slide 2 of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the table shape ("Tabla 1") on the slide.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $tableShape = $shape
        break
    }
}

$tbl = $tableShape.Table

# Append a new row at the bottom of the table and fill in its cells.
$newRow = $tbl.Rows.Add()
$rowIndex = $tbl.Rows.Count

$tbl.Cell($rowIndex, 1).Shape.TextFrame.TextRange.Text = "Febrero"
$tbl.Cell($rowIndex, 2).Shape.TextFrame.TextRange.Text = "27-02-15"
$tbl.Cell($rowIndex, 3).Shape.TextFrame.TextRange.Text = "27-02-15"
